# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data refresh to the Gungnir_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 3704042.8
$ws.Range("I92").Value = 9259510
$ws.Range("J92").Value = 398.33334
$ws.Range("K92").Value = 9259510
$ws.Range("L92").Value = 398.33334
$ws.Range("M92").Value = -9258262
$ws.Range("N92").Value = -2894.33334
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""
$ws.Range("H95").Value = 19999
$ws.Range("J95").Value = 19999
$ws.Range("L95").Value = 19999
$ws.Range("N95").Value = -25491
$ws.Range("H96").Value = 1968.7142
$ws.Range("I96").Value = 2250.6667
$ws.Range("J96").Value = 1757.25
$ws.Range("K96").Value = 6752.000100000001
$ws.Range("L96").Value = 5271.75
$ws.Range("M96").Value = -5379.000100000001
$ws.Range("N96").Value = -8017.75
$ws.Range("H97").Value = 251252270
$ws.Range("J97").Value = 251252270
$ws.Range("L97").Value = 753756810
$ws.Range("N97").Value = -753757802
$ws.Range("H99").Value = 402.63635
$ws.Range("I99").Value = 364.2857
$ws.Range("J99").Value = 469.75
$ws.Range("K99").Value = 1092.8571
$ws.Range("L99").Value = 1409.25
$ws.Range("M99").Value = 405.1428999999998
$ws.Range("N99").Value = -4405.25
$ws.Range("H100").Value = 7886.316
$ws.Range("I100").Value = 10217.538
$ws.Range("K100").Value = 10217.538
$ws.Range("M100").Value = -9676.538
$ws.Range("H121").Value = 872.7
$ws.Range("I121").Value = 340.66666
$ws.Range("J121").Value = 1100.7142
$ws.Range("K121").Value = 1021.99998
$ws.Range("L121").Value = 3302.1426
$ws.Range("M121").Value = 725.0000200000001
$ws.Range("N121").Value = -6796.142599999999
$ws.Range("H131").Value = 4619.857
$ws.Range("I131").Value = 837
$ws.Range("J131").Value = 6721.4443
$ws.Range("K131").Value = 2511
$ws.Range("L131").Value = 20164.3329
$ws.Range("M131").Value = 2529
$ws.Range("N131").Value = -30244.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1500
$ws.Range("I35").Value = 1500
$ws.Range("K35").Value = 1500
$ws.Range("M35").Value = -1094
$ws.Range("H86").Value = 1000000
$ws.Range("J86").Value = 1000000
$ws.Range("L86").Value = 1000000
$ws.Range("N86").Value = -1002372
$ws.Range("H89").Value = 1000000
$ws.Range("J89").Value = 1000000
$ws.Range("L89").Value = 3000000
$ws.Range("N89").Value = -3011856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3371177.5
$ws.Range("I134").Value = 1035.5
$ws.Range("J134").Value = 12358223
$ws.Range("K134").Value = 3106.5
$ws.Range("L134").Value = 37074669
$ws.Range("M134").Value = -571.5
$ws.Range("N134").Value = -37079739

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2141.5
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 2772
$ws.Range("K22").Value = 250
$ws.Range("L22").Value = 2772
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = -3472

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 16692640
$ws.Range("J22").Value = 31168
$ws.Range("L22").Value = 93504
$ws.Range("N22").Value = -93842
$ws.Range("H27").Value = 16692640
$ws.Range("J27").Value = 31168
$ws.Range("L27").Value = 93504
$ws.Range("N27").Value = -93708
$ws.Range("H58").Value = 65876
$ws.Range("J58").Value = 82095
$ws.Range("L58").Value = 246285
$ws.Range("N58").Value = -246541
$ws.Range("H92").Value = 8310
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 8310
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 24930
$ws.Range("M92").Value = ""
$ws.Range("N92").Value = -27426
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = ""
$ws.Range("H96").Value = 8666.666999999999
$ws.Range("J96").Value = 8666.666999999999
$ws.Range("L96").Value = 26000.001
$ws.Range("N96").Value = -30118.001
$ws.Range("H98").Value = 916.6667
$ws.Range("I98").Value = 750
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 2250
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -752
$ws.Range("N98").Value = -5996
$ws.Range("H99").Value = 1450
$ws.Range("I99").Value = 900
$ws.Range("K99").Value = 2700
$ws.Range("M99").Value = -454
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""
$ws.Range("H106").Value = 7142.857
$ws.Range("J106").Value = 7400
$ws.Range("L106").Value = 22200
$ws.Range("N106").Value = -24092
$ws.Range("H121").Value = 786.6
$ws.Range("I121").Value = 400
$ws.Range("J121").Value = 883.25
$ws.Range("K121").Value = 1200
$ws.Range("L121").Value = 2649.75
$ws.Range("M121").Value = 110
$ws.Range("N121").Value = -5269.75
$ws.Range("H131").Value = 799.37
$ws.Range("J131").Value = 806.01044
$ws.Range("L131").Value = 2418.03132
$ws.Range("N131").Value = -12498.03132
$ws.Range("H140").Value = 10001524
$ws.Range("I140").Value = 13158952
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 39476856
$ws.Range("L140").Value = 9000
$ws.Range("M140").Value = -39471676
$ws.Range("N140").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 30000
$ws.Range("J88").Value = 30000
$ws.Range("L88").Value = 30000
$ws.Range("N88").Value = -30902
$ws.Range("H91").Value = 30000
$ws.Range("J91").Value = 30000
$ws.Range("L91").Value = 30000
$ws.Range("N91").Value = -33120
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 983.6
$ws.Range("I16").Value = 979.625
$ws.Range("J16").Value = 999.5
$ws.Range("K16").Value = 979.625
$ws.Range("L16").Value = 999.5
$ws.Range("M16").Value = -809.625
$ws.Range("N16").Value = -1339.5
$ws.Range("H87").Value = 2000
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = ""
$ws.Range("H90").Value = 2000
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = ""
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = ""
$ws.Range("H93").Value = 1452.7858
$ws.Range("I93").Value = 1148.4286
$ws.Range("J93").Value = 1757.1428
$ws.Range("K93").Value = 1148.4286
$ws.Range("L93").Value = 1757.1428
$ws.Range("M93").Value = 99.57140000000004
$ws.Range("N93").Value = -4253.1428
$ws.Range("H95").Value = 35175
$ws.Range("J95").Value = 35175
$ws.Range("L95").Value = 35175
$ws.Range("M95").Value = -40667
$ws.Range("H97").Value = 14806.125
$ws.Range("J97").Value = 14806.125
$ws.Range("L97").Value = 14806.125
$ws.Range("N97").Value = -16788.125
$ws.Range("H98").Value = 40000
$ws.Range("J98").Value = 40000
$ws.Range("L98").Value = 40000
$ws.Range("N98").Value = -45990
$ws.Range("H99").Value = 28000
$ws.Range("J99").Value = 28000
$ws.Range("L99").Value = 28000
$ws.Range("N99").Value = -33990
$ws.Range("H100").Value = 2522.1875
$ws.Range("I100").Value = 2933.8333
$ws.Range("J100").Value = 2427.1924
$ws.Range("K100").Value = 2933.8333
$ws.Range("L100").Value = 2427.1924
$ws.Range("M100").Value = -2392.8333
$ws.Range("N100").Value = -3509.1924
$ws.Range("H101").Value = 21675.5
$ws.Range("J101").Value = 21675.5
$ws.Range("L101").Value = 21675.5
$ws.Range("N101").Value = -28165.5
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = ""
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 49750
$ws.Range("J92").Value = 49750
$ws.Range("L92").Value = 49750
$ws.Range("N92").Value = -54742
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H96").Value = 1674
$ws.Range("I96").Value = 1440
$ws.Range("J96").Value = 1908
$ws.Range("K96").Value = 1440
$ws.Range("L96").Value = 1908
$ws.Range("M96").Value = -67
$ws.Range("N96").Value = -4654
$ws.Range("H100").Value = 521.25
$ws.Range("I100").Value = 420
$ws.Range("K100").Value = 840
$ws.Range("M100").Value = -299
